$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48; existing rows 48-104 shift down to 49-105
$ws.Rows.Item(48).Insert()

# Fill in the new row 48 with its data
$ws.Cells.Item(48, 1).Value2 = 11
$ws.Cells.Item(48, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(48, 3).Value2 = "Bíobío"
$ws.Cells.Item(48, 4).Value2 = 44601
$ws.Cells.Item(48, 5).Value2 = 8
$ws.Cells.Item(48, 6).Value2 = "Fruta"
$ws.Cells.Item(48, 7).Value2 = 100109
$ws.Cells.Item(48, 8).Value2 = "Uva"
$ws.Cells.Item(48, 9).Value2 = 100109001
$ws.Cells.Item(48, 10).Value2 = "Uva"
$ws.Cells.Item(48, 11).Value2 = "Superior Seedless"
$ws.Cells.Item(48, 12).Value2 = "Primera"
$ws.Cells.Item(48, 13).Value2 = 200
$ws.Cells.Item(48, 14).Value2 = 9000
$ws.Cells.Item(48, 15).Value2 = 10000
$ws.Cells.Item(48, 16).Value2 = 9500
$ws.Cells.Item(48, 17).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(48, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(48, 19).Value2 = 528
$ws.Cells.Item(48, 20).Value2 = 18

Write-Host ("Dimension used range rows: " + $ws.UsedRange.Rows.Count)
